$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 54,6

$arr[0,0] = 0.003208871385164791
$arr[0,1] = 0.002571899574220771
$arr[0,2] = 261.3203778131603
$arr[0,3] = 10.19245300693656
$arr[0,4] = 0
$arr[0,5] = 271.5186115910562

$arr[1,0] = 3.286832544864788
$arr[1,1] = 1.655778082260271
$arr[1,2] = 0.7527432677738641
$arr[1,3] = 0.4942365360607697
$arr[1,4] = 1
$arr[1,5] = 6.189590430959694

$arr[2,0] = 3.286832544864788
$arr[2,1] = 1.655778082260271
$arr[2,2] = 0.7527432677738641
$arr[2,3] = 0.4942365360607697
$arr[2,4] = 1
$arr[2,5] = 6.189590430959694

$arr[3,0] = 3.286832544864788
$arr[3,1] = 1.655778082260271
$arr[3,2] = 22.3905356188092
$arr[3,3] = 0.4942365360607697
$arr[3,4] = 0
$arr[3,5] = 27.82738278199502

$arr[4,0] = 1.455362044514542
$arr[4,1] = 0.306821227259698
$arr[4,2] = 22.3905356188092
$arr[4,3] = 10.19245300693656
$arr[4,4] = 1
$arr[4,5] = 34.34517189751999

$arr[5,0] = 1.455362044514542
$arr[5,1] = 0.306821227259698
$arr[5,2] = 0.1494219747398047
$arr[5,3] = 10.19245300693656
$arr[5,4] = 0
$arr[5,5] = 12.1040582534506

$arr[6,0] = 0.1190320826869504
$arr[6,1] = 0.04071648406533734
$arr[6,2] = 0.1494219747398047
$arr[6,3] = 0.4942365360607697
$arr[6,4] = 0
$arr[6,5] = 0.8034070775528621

$arr[7,0] = 3.286832544864788
$arr[7,1] = 1.655778082260271
$arr[7,2] = 0.1494219747398047
$arr[7,3] = 0.4942365360607697
$arr[7,4] = 1
$arr[7,5] = 5.586269137925634

$arr[8,0] = 3.286832544864788
$arr[8,1] = 1.655778082260271
$arr[8,2] = 0.7527432677738641
$arr[8,3] = 0.4942365360607697
$arr[8,4] = 1
$arr[8,5] = 6.189590430959694

$arr[9,0] = 1.455362044514542
$arr[9,1] = 1.655778082260271
$arr[9,2] = 3.537761648806719
$arr[9,3] = 0.4942365360607697
$arr[9,4] = 1
$arr[9,5] = 7.143138311642302

$arr[10,0] = 3.286832544864788
$arr[10,1] = 1.655778082260271
$arr[10,2] = 3.537761648806719
$arr[10,3] = 0.4942365360607697
$arr[10,4] = 1
$arr[10,5] = 8.974608811992548

$arr[11,0] = 0.6606524410359556
$arr[11,1] = 1.655778082260271
$arr[11,2] = 22.3905356188092
$arr[11,3] = 0.4942365360607697
$arr[11,4] = 1
$arr[11,5] = 25.20120267816619

$arr[12,0] = 3.286832544864788
$arr[12,1] = 1.655778082260271
$arr[12,2] = 3.537761648806719
$arr[12,3] = 0.4942365360607697
$arr[12,4] = 1
$arr[12,5] = 8.974608811992548

$arr[13,0] = 3.286832544864788
$arr[13,1] = 1.655778082260271
$arr[13,2] = 0.1494219747398047
$arr[13,3] = 10.19245300693656
$arr[13,4] = 1
$arr[13,5] = 15.28448560880142

$arr[14,0] = 1.455362044514542
$arr[14,1] = 1.655778082260271
$arr[14,2] = 0.7527432677738641
$arr[14,3] = 0.4942365360607697
$arr[14,4] = 0
$arr[14,5] = 4.358119930609447

$arr[15,0] = 3.286832544864788
$arr[15,1] = 1.655778082260271
$arr[15,2] = 0.7527432677738641
$arr[15,3] = 0.4942365360607697
$arr[15,4] = 1
$arr[15,5] = 6.189590430959694

$arr[16,0] = 3.286832544864788
$arr[16,1] = 1.655778082260271
$arr[16,2] = 0.7527432677738641
$arr[16,3] = 0.4942365360607697
$arr[16,4] = 1
$arr[16,5] = 6.189590430959694

$arr[17,0] = 1.455362044514542
$arr[17,1] = 1.655778082260271
$arr[17,2] = 3.537761648806719
$arr[17,3] = 10.19245300693656
$arr[17,4] = 1
$arr[17,5] = 16.84135478251809

$arr[18,0] = 0.1190320826869504
$arr[18,1] = 0.306821227259698
$arr[18,2] = 0.7527432677738641
$arr[18,3] = 0.4942365360607697
$arr[18,4] = 1
$arr[18,5] = 1.672833113781282

$arr[19,0] = 3.286832544864788
$arr[19,1] = 1.655778082260271
$arr[19,2] = 3.537761648806719
$arr[19,3] = 0.4942365360607697
$arr[19,4] = 1
$arr[19,5] = 8.974608811992548

$arr[20,0] = 1.455362044514542
$arr[20,1] = 1.655778082260271
$arr[20,2] = 0.7527432677738641
$arr[20,3] = 0.4942365360607697
$arr[20,4] = 1
$arr[20,5] = 4.358119930609447

$arr[21,0] = 3.286832544864788
$arr[21,1] = 10.34677158129881
$arr[21,2] = 0.1494219747398047
$arr[21,3] = 10.19245300693656
$arr[21,4] = 1
$arr[21,5] = 23.97547910783996

$arr[22,0] = 3.286832544864788
$arr[22,1] = 10.34677158129881
$arr[22,2] = 0.7527432677738641
$arr[22,3] = 10.19245300693656
$arr[22,4] = 1
$arr[22,5] = 24.57880040087402

$arr[23,0] = 3.286832544864788
$arr[23,1] = 1.655778082260271
$arr[23,2] = 3.537761648806719
$arr[23,3] = 0.4942365360607697
$arr[23,4] = 1
$arr[23,5] = 8.974608811992548

$arr[24,0] = 1.455362044514542
$arr[24,1] = 1.655778082260271
$arr[24,2] = 0.7527432677738641
$arr[24,3] = 0.4942365360607697
$arr[24,4] = 1
$arr[24,5] = 4.358119930609447

$arr[25,0] = 3.286832544864788
$arr[25,1] = 1.655778082260271
$arr[25,2] = 3.537761648806719
$arr[25,3] = 0.4942365360607697
$arr[25,4] = 1
$arr[25,5] = 8.974608811992548

$arr[26,0] = 1.455362044514542
$arr[26,1] = 1.655778082260271
$arr[26,2] = 0.7527432677738641
$arr[26,3] = 0.4942365360607697
$arr[26,4] = 0
$arr[26,5] = 4.358119930609447

$arr[27,0] = 3.286832544864788
$arr[27,1] = 1.655778082260271
$arr[27,2] = 0.7527432677738641
$arr[27,3] = 0.4942365360607697
$arr[27,4] = 0
$arr[27,5] = 6.189590430959694

$arr[28,0] = 0.1190320826869504
$arr[28,1] = 0.04071648406533734
$arr[28,2] = 0.1494219747398047
$arr[28,3] = 0.4942365360607697
$arr[28,4] = 0
$arr[28,5] = 0.8034070775528621

$arr[29,0] = 1.455362044514542
$arr[29,1] = 1.655778082260271
$arr[29,2] = 3.537761648806719
$arr[29,3] = 10.19245300693656
$arr[29,4] = 1
$arr[29,5] = 16.84135478251809

$arr[30,0] = 3.286832544864788
$arr[30,1] = 1.655778082260271
$arr[30,2] = 0.7527432677738641
$arr[30,3] = 0.4942365360607697
$arr[30,4] = 0
$arr[30,5] = 6.189590430959694

$arr[31,0] = 3.286832544864788
$arr[31,1] = 1.655778082260271
$arr[31,2] = 0.1494219747398047
$arr[31,3] = 0.4942365360607697
$arr[31,4] = 1
$arr[31,5] = 5.586269137925634

$arr[32,0] = 3.286832544864788
$arr[32,1] = 1.655778082260271
$arr[32,2] = 0.7527432677738641
$arr[32,3] = 0.4942365360607697
$arr[32,4] = 1
$arr[32,5] = 6.189590430959694

$arr[33,0] = 0.6606524410359556
$arr[33,1] = 1.655778082260271
$arr[33,2] = 0.7527432677738641
$arr[33,3] = 0.4942365360607697
$arr[33,4] = 1
$arr[33,5] = 3.56341032713086

$arr[34,0] = 1.455362044514542
$arr[34,1] = 0.306821227259698
$arr[34,2] = 0.1494219747398047
$arr[34,3] = 0.4942365360607697
$arr[34,4] = 1
$arr[34,5] = 2.405841782574814

$arr[35,0] = 3.286832544864788
$arr[35,1] = 1.655778082260271
$arr[35,2] = 22.3905356188092
$arr[35,3] = 0.4942365360607697
$arr[35,4] = 1
$arr[35,5] = 27.82738278199502

$arr[36,0] = 0.2917716402565462
$arr[36,1] = 0.306821227259698
$arr[36,2] = 0.7527432677738641
$arr[36,3] = 0.4942365360607697
$arr[36,4] = 1
$arr[36,5] = 1.845572671350878

$arr[37,0] = 3.286832544864788
$arr[37,1] = 1.655778082260271
$arr[37,2] = 22.3905356188092
$arr[37,3] = 0.4942365360607697
$arr[37,4] = 1
$arr[37,5] = 27.82738278199502

$arr[38,0] = 3.286832544864788
$arr[38,1] = 1.655778082260271
$arr[38,2] = 0.1494219747398047
$arr[38,3] = 0.4942365360607697
$arr[38,4] = 1
$arr[38,5] = 5.586269137925634

$arr[39,0] = 1.455362044514542
$arr[39,1] = 1.655778082260271
$arr[39,2] = 0.7527432677738641
$arr[39,3] = 0.4942365360607697
$arr[39,4] = 1
$arr[39,5] = 4.358119930609447

$arr[40,0] = 0.1190320826869504
$arr[40,1] = 1.655778082260271
$arr[40,2] = 0.7527432677738641
$arr[40,3] = 0.4942365360607697
$arr[40,4] = 1
$arr[40,5] = 3.021789968781855

$arr[41,0] = 0.6606524410359556
$arr[41,1] = 1.655778082260271
$arr[41,2] = 0.7527432677738641
$arr[41,3] = 0.4942365360607697
$arr[41,4] = 1
$arr[41,5] = 3.56341032713086

$arr[42,0] = 0.6606524410359556
$arr[42,1] = 1.655778082260271
$arr[42,2] = 0.7527432677738641
$arr[42,3] = 0.4942365360607697
$arr[42,4] = 1
$arr[42,5] = 3.56341032713086

$arr[43,0] = 3.286832544864788
$arr[43,1] = 1.655778082260271
$arr[43,2] = 0.7527432677738641
$arr[43,3] = 0.4942365360607697
$arr[43,4] = 1
$arr[43,5] = 6.189590430959694

$arr[44,0] = 1.455362044514542
$arr[44,1] = 1.655778082260271
$arr[44,2] = 0.7527432677738641
$arr[44,3] = 0.4942365360607697
$arr[44,4] = 1
$arr[44,5] = 4.358119930609447

$arr[45,0] = 0.6606524410359556
$arr[45,1] = 1.655778082260271
$arr[45,2] = 22.3905356188092
$arr[45,3] = 0.4942365360607697
$arr[45,4] = 1
$arr[45,5] = 25.20120267816619

$arr[46,0] = 0.6606524410359556
$arr[46,1] = 0.04071648406533734
$arr[46,2] = 0.1494219747398047
$arr[46,3] = 0.4942365360607697
$arr[46,4] = 1
$arr[46,5] = 1.345027435901867

$arr[47,0] = 3.286832544864788
$arr[47,1] = 1.655778082260271
$arr[47,2] = 0.1494219747398047
$arr[47,3] = 0.4942365360607697
$arr[47,4] = 1
$arr[47,5] = 5.586269137925634

$arr[48,0] = 3.286832544864788
$arr[48,1] = 1.655778082260271
$arr[48,2] = 3.537761648806719
$arr[48,3] = 0.4942365360607697
$arr[48,4] = 1
$arr[48,5] = 8.974608811992548

$arr[49,0] = 1.455362044514542
$arr[49,1] = 1.655778082260271
$arr[49,2] = 0.7527432677738641
$arr[49,3] = 0.4942365360607697
$arr[49,4] = 1
$arr[49,5] = 4.358119930609447

$arr[50,0] = 0.6606524410359556
$arr[50,1] = 1.655778082260271
$arr[50,2] = 0.7527432677738641
$arr[50,3] = 0.4942365360607697
$arr[50,4] = 1
$arr[50,5] = 3.56341032713086

$arr[51,0] = 3.286832544864788
$arr[51,1] = 1.655778082260271
$arr[51,2] = 0.1494219747398047
$arr[51,3] = 0.4942365360607697
$arr[51,4] = 1
$arr[51,5] = 5.586269137925634

$arr[52,0] = 3.286832544864788
$arr[52,1] = 1.655778082260271
$arr[52,2] = 3.537761648806719
$arr[52,3] = 0.4942365360607697
$arr[52,4] = 1
$arr[52,5] = 8.974608811992548

$arr[53,0] = 0.6606524410359556
$arr[53,1] = 1.655778082260271
$arr[53,2] = 0.7527432677738641
$arr[53,3] = 0.4942365360607697
$arr[53,4] = 1
$arr[53,5] = 3.56341032713086

$ws.Range("B2:G55").Value = $arr

Write-Output "applied updates"
